$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a brand-new row above the current row 5, pushing the existing
# rows 5-7 down to rows 6-8 (values, formatting and style all move with them).
$ws.Rows.Item(5).Insert()

# --- Rows 3 & 4: clear explicit cell style (falls back to default "Normal")
# and change the Curation status (column S) from "Proposed" to "External".
$ws.Range("A3:V4").Style = "Normal"
$ws.Range("S3").Value = "External"
$ws.Range("S4").Value = "External"

# --- Populate the newly inserted row 5 with the new "data transformation" entry.
$ws.Range("A5").Value = "OBI:0200000"
$ws.Range("B5").Value = "data transformation"
$ws.Range("C5").Value = "A planned process that produces output data from input data."
$ws.Range("D5").Value = "planned process [COB:0000082]"
$ws.Range("S5").Value = "Proposed"
$ws.Range("V5").Value = "BG"
